# Disaggregation of commodity Copper
#
# 1) Rename the shared-string label "Copper ores and concentrates" -> "Copper".
#    That label is used (cell C4) on every year-tab (2000..2100), so every
#    occurrence must be updated for the shared string table to collapse back
#    down to a single "Copper" entry.
# 2) A handful of tabs also carry a replacement cached value in D4 (last-digit
#    float re-computation that came along with the relabeling).

$wb = $excel.ActiveWorkbook

# --- 1) Update the commodity label on every sheet -------------------------
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("C4").Value = "Copper"
}

# --- 2) Patch the handful of recomputed D4 cached values -------------------
$updatedValues = @{
    2021 = 30391.99057451256
    2023 = 50021.9813159201
    2025 = 57877.16490272013
    2028 = 84055.74857747872
    2041 = 455285.6495998815
    2044 = 907293.3796566341
    2074 = 1771147.898692237
    2090 = 1847140.799864977
    2092 = 1876913.119977531
}

foreach ($sheetName in $updatedValues.Keys) {
    $ws = $wb.Worksheets.Item("$sheetName")
    $ws.Range("D4").Value = $updatedValues[$sheetName]
}
